$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Duplicate the formatting of row 34/row 33's pattern down onto a brand new
#    row 35 so the new row inherits the correct alternating-style banding
#    (row 33 uses the "odd" style band: s=7 normally, s=8 for M/O/P/AC).
# ---------------------------------------------------------------------------
$ws.Range("A33:AK33").Copy($ws.Range("A35:AK35"))

# ---------------------------------------------------------------------------
# 2) Fill in the real values for the newly added service record (row 35).
# ---------------------------------------------------------------------------
$ws.Range("A35").Value = 33
$ws.Range("B35").Value = "服務"
$ws.Range("C35").Value = 2025071115
$ws.Range("D35").Value = ""
$ws.Range("E35").Value = ""
$ws.Range("F35").Value = 3890
$ws.Range("G35").Value = "三重美堤店"
$ws.Range("H35").Value = "新北市三重區"
$ws.Range("I35").Value = ""
$ws.Range("J35").Value = ""
$ws.Range("K35").Value = ""
$ws.Range("L35").Value = ""
$ws.Range("M35").Value = ""
$ws.Range("N35").Value = ""
$ws.Range("O35").Value = ""
$ws.Range("P35").Value = ""
$ws.Range("Q35").Value = "THILF03890"
$ws.Range("R35").Value = "新北一"
$ws.Range("S35").Value = "吳宗鴻"
$ws.Range("T35").Value = 1
$ws.Range("U35").Value = "已完工"
$ws.Range("V35").Value = "2025-07-07 17:31:46"
$ws.Range("W35").Value = "2025-07-07 14:00:00"
$ws.Range("X35").Value = "2025-07-07 17:30:00"
$ws.Range("Y35").Value = ""
$ws.Range("Z35").Value = 3.5
$ws.Range("AA35").Value = ""
$ws.Range("AB35").Value = "到場處理"
$ws.Range("AC35").Value = "重新回裝已完工"
$ws.Range("AD35").Value = ""
$ws.Range("AE35").Value = "O"
$ws.Range("AF35").Value = ""
$ws.Range("AG35").Value = ""
$ws.Range("AH35").Value = ""
$ws.Range("AI35").Value = ""
$ws.Range("AJ35").Value = ""
$ws.Range("AK35").Value = "O"

# ---------------------------------------------------------------------------
# 3) Row 34's "work content" (P34/AC34) now needs to wrap its text, same as
#    other rows in the sheet that already wrap the notes column.
# ---------------------------------------------------------------------------
$ws.Range("P34").WrapText = $true
$ws.Range("AC34").WrapText = $true

# ---------------------------------------------------------------------------
# 4) Update Print_Area (localSheetId 0 -> 'Report' sheet) to include new row 35.
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Report!Print_Area") {
        $n.RefersTo = "='Report'!`$A`$1:`$AK`$35"
    }
}

# ---------------------------------------------------------------------------
# 5) Move the active selection to the new last row, like the saved workbook.
# ---------------------------------------------------------------------------
$ws.Range("A35").Select()
